# Auto-update draw results: append the 2025-12-08 Pick 4 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 83

# The sheet stores every column as plain text (dates, phase codes, and the
# result string all look numeric/date-like to Excel's auto-detection), so
# force text formatting on the new row before writing values to avoid them
# being reinterpreted as dates/numbers.
$rowRange = $ws.Range("A" + $newRow + ":E" + $newRow)
$rowRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-12-08"
$ws.Cells.Item($newRow, 2).Value = "Pick 4"
$ws.Cells.Item($newRow, 3).Value = "251208"
$ws.Cells.Item($newRow, 4).Value = "4-5-8-3"
$ws.Cells.Item($newRow, 5).Value = "2025-12-08T21:42:08.107+04:00"

# Keep the sheet's "numbers stored as text" warning suppressed over the
# now-larger used range (best effort - harmless if the host doesn't wire
# this error-checking toggle through).
try {
    $ws.Range("A1:E" + $newRow).Errors.Item(9).Ignore = $true
} catch {
}
